# Generate Report for Handoff
# Moves the localization status from "In Translation" to "Ready for handoff"
# and refreshes the associated timestamps, widening the "Status" columns
# to fit the new (longer) status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-17 18:57:31"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-17 18:57:26"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-17 18:57:31"

# --- Widen the "Status" columns so the longer text fits ---------------------
# (target stored width ~17.216 characters; ColumnWidth is expressed
# without the standard 5-pixel cell padding that Excel adds on save)
$newColumnWidth = 17.2159881591797 - 0.8333333333333334

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
